$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '73.019.53'
$ws.Range("E2").Value = '  +1.25%  '
$ws.Range("D3").Value = '3.960.24'
$ws.Range("E3").Value = '  -1.89%  '
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '614.58'
$ws.Range("E5").Value = '  +13.74%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '166.29'
$ws.Range("E6").Value = '  +9.61%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.681'
$ws.Range("E7").Value = '  -2.70%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.754'
$ws.Range("E9").Value = '  -0.03%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.184'
$ws.Range("E10").Value = '  +6.92%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '56.07'
$ws.Range("E11").Value = '  +4.13%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000333'
$ws.Range("E12").Value = '  +1.45%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '11.12'
$ws.Range("E13").Value = '  +2.11%  '
$ws.Range("D14").Value = '4.593.09'
$ws.Range("E14").Value = '  -1.93%  '
$ws.Range("D15").Value = '3.964.67'
$ws.Range("E15").Value = '  -2.00%  '
$ws.Range("E16").Value = '  +3.72%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.10'
$ws.Range("E17").Value = '  -1.96%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '20.46'
$ws.Range("E18").Value = '  -0.69%  '
$ws.Range("B19").Value = 'WrappedBTC'
$ws.Range("C19").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D19").Value = '72.859.51'
$ws.Range("E19").Value = '  +1.04%  '
$ws.Range("B20").Value = 'TRON'
$ws.Range("C20").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.131'
$ws.Range("E20").Value = '  -0.43%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '438.68'
$ws.Range("E21").Value = '  -2.02%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.91'
$ws.Range("E22").Value = '  +15.60%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '95.86'
$ws.Range("E23").Value = '  -1.73%  '
$ws.Range("E24").Value = '  -3.58%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '14.13'
$ws.Range("E25").Value = '  -3.49%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '4.07'
$ws.Range("E26").Value = '  -5.89%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.07'
$ws.Range("E27").Value = '  -2.22%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.97'
$ws.Range("E28").Value = '  +0.12%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.50'
$ws.Range("E29").Value = '  -2.77%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '35.98'
$ws.Range("E30").Value = '  -3.13%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.00'
$ws.Range("E31").Value = '  -1.96%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '13.63'
$ws.Range("E32").Value = '  +0.21%  '
$ws.Range("E33").Value = '  +19.26%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.130'
$ws.Range("E34").Value = '  -3.75%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '47.49'
$ws.Range("E35").Value = '  -3.49%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '70.33'
$ws.Range("E36").Value = '  +5.23%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '649.76'
$ws.Range("E37").Value = '  -4.30%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.430'
$ws.Range("E38").Value = '  -5.69%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.46'
$ws.Range("E39").Value = '  +2.29%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.00'
$ws.Range("E40").Value = '  -0.13%  '
$ws.Range("E41").Value = '  -2.21%  '
$ws.Range("E42").Value = '  +0.02%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0483'
$ws.Range("E43").Value = '  -2.17%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '10.66'
$ws.Range("E44").Value = '  -5.00%  '
$ws.Range("E45").Value = '  -5.66%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.149'
$ws.Range("E46").Value = '  -1.58%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.43'
$ws.Range("E47").Value = '  +3.44%  '
$ws.Range("B48").Value = 'dogwifhat'
$ws.Range("C48").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.88'
$ws.Range("E48").Value = '  +25.33%  '
$ws.Range("B49").Value = 'Fetch.AI'
$ws.Range("C49").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.60'
$ws.Range("E49").Value = '  -1.41%  '
$ws.Range("D50").Value = '2.833.55'
$ws.Range("E50").Value = '  +3.21%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '149.66'
$ws.Range("E51").Value = '  +1.76%  '
